# gerar conteúdo por módulo de análise
#
# The underlying analysis data was regenerated; for the UFs that are tied
# (zero / missing values in a given module) the resulting ranking order of
# MS / MT / TO / SE / AC / RO shifts slightly in a few of the per-module
# sheets. Re-apply the new label for each affected ranking row.

$wb = $excel.ActiveWorkbook

# tot-arrecad
$ws = $wb.Worksheets.Item(2)
$ws.Range("A20").Value = "MS"
$ws.Range("A21").Value = "MT"
$ws.Range("A23").Value = "TO"

# avg-arrecad
$ws = $wb.Worksheets.Item(3)
$ws.Range("A2").Value = "MS"
$ws.Range("A21").Value = "MT"
$ws.Range("A23").Value = "TO"
$ws.Range("A25").Value = "AC"
$ws.Range("A26").Value = "RO"

# max-arrecad
$ws = $wb.Worksheets.Item(4)
$ws.Range("A14").Value = "MS"
$ws.Range("A21").Value = "MT"
$ws.Range("A23").Value = "AC"
$ws.Range("A24").Value = "TO"
$ws.Range("A25").Value = "SE"
$ws.Range("A26").Value = "RO"

# tx-sucesso
$ws = $wb.Worksheets.Item(5)
$ws.Range("A2").Value = "MT"
$ws.Range("A15").Value = "MS"
$ws.Range("A23").Value = "TO"
